# 989: Adds CMS to test imports and cms table to extract process
$wb = $excel.ActiveWorkbook

# --- Update selection on WMT_Extract (sheet1): drop the frozen/top-left
# scroll position and move the selection to C24 ---
$wsExtract = $wb.Worksheets.Item("WMT_Extract")
[void]$wsExtract.Range("C24").Select()

# --- Add the new CMS worksheet as the last tab ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCms = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$wsCms.Name = "CMS"

# Header row values (new shared strings 52-65)
$headers = @(
    "Contact_ID`n",
    "Contact_Date`n",
    "Contact_Type_Code",
    "Contact_Type_Desc",
    "Contact_Staff_Name",
    "Contact_Staff_Key",
    "Contact_Staff_Grade",
    "Contact_Team_Key",
    "Contact_Provider_Code",
    "OM_Name`n",
    "OM_Key`n",
    "OM_Grade`n",
    "OM_Team_Key`n",
    "OM_Provider_Code`n"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsCms.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Style the header row: 9pt Arial dark-grey text, white fill, left aligned,
# with a taller custom row height.
$headerRange = $wsCms.Range("A1:N1")
$headerRange.Font.Name = "Arial"
$headerRange.Font.Size = 9
$headerRange.Font.Color = 0x333333
$headerRange.Interior.Color = 0xFFFFFF
$headerRange.Interior.PatternColor = 0xFFFFFF
$headerRange.HorizontalAlignment = -4131
$wsCms.Rows.Item(1).RowHeight = 23.25

# Select C4 on the new sheet so it ends up as the active tab/selection
[void]$wsCms.Range("C4").Select()
